# Append a new row (date + gold-price text) at the bottom of the data table,
# mirroring the style of the immediately preceding row, while writing the
# cell content as plain text (not an auto-converted date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$newDate = "12-12-2025"
$newText = "The price of gold in India today is ₹13,266 per gram for 24 karat gold, ₹12,160 per gram for 22 karat gold and ₹9,949 per gram for 18 karat gold (also called 999 gold)."

$srcA = $ws.Cells.Item($lastRow, 1)
$srcB = $ws.Cells.Item($lastRow, 2)
$dstA = $ws.Cells.Item($newRow, 1)
$dstB = $ws.Cells.Item($newRow, 2)

# Clone the formatting (border/wrap styles) of the row above onto the new row.
$srcA.Copy($dstA) | Out-Null
$srcB.Copy($dstB) | Out-Null

# Stage the literal text values on scratch cells far away via a text formula
# result (so Excel's "looks like a date" auto-conversion never kicks in),
# then paste only the values onto the freshly-formatted destination cells.
$scratchA = $ws.Cells.Item(2000, 1)
$scratchB = $ws.Cells.Item(2000, 2)

$scratchA.Formula = '="' + $newDate + '"'
$scratchB.Formula = '="' + $newText + '"'

$scratchA.Copy() | Out-Null
$dstA.PasteSpecial(-4163) | Out-Null

$scratchB.Copy() | Out-Null
$dstB.PasteSpecial(-4163) | Out-Null

$ws.Rows.Item(2000).Delete() | Out-Null
